$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2-92) lists one row per test script (DEC_xxxx).
# We're appending scripts DEC_0176 .. DEC_0190 (15 new rows) to the
# "Tests_AdmInstituciones" block, right before the trailing
# blank-line / ASDF-Verity1.0 / USUARIO-PASSWORD footer rows that sit
# at the bottom of the sheet (previously rows 93-96).
#
# Rows 93 and 94 are already blank spacer rows, so we reuse them as the
# first two new data rows and insert 13 more rows above the footer to
# fit the rest, which pushes the footer down from rows 95-96 to 109-110
# and leaves a single blank spacer row (108) in between, matching the
# target layout.
$ws.Range("A95:A108").EntireRow.Insert()

$codes = @(
    "DEC_0176", "DEC_0177", "DEC_0178", "DEC_0179", "DEC_0180",
    "DEC_0181", "DEC_0182", "DEC_0183", "DEC_0184", "DEC_0185",
    "DEC_0186", "DEC_0187", "DEC_0188", "DEC_0189", "DEC_0190"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = 93 + $i
    $ws.Cells.Item($r, 1).Value = $codes[$i]
    $ws.Cells.Item($r, 2).Value = "13712759-8"
    $ws.Cells.Item($r, 3).Value = "Verity1.1"
    $ws.Cells.Item($r, 4).Value = "SIN_DATO"
    $ws.Cells.Item($r, 5).Value = "SIN_DATO"
    $ws.Cells.Item($r, 6).Value = "SIN_DATO"
    $ws.Cells.Item($r, 7).Value = "SIN_DATO"
    $ws.Cells.Item($r, 8).Value = "SIN_DATO"
    $ws.Cells.Item($r, 9).Value = "SIN_DATO"
    $ws.Cells.Item($r, 10).Value = "SIN_DATO"
}

# Update the view to match where the editor ended up after adding the rows.
$ws.Range("H101").Select()
